$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as literal text
# (matching the source inlineStr cells) instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.416.42"
$ws.Range("E2").Value = "  +8.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.157.94"
$ws.Range("E3").Value = "  +6.03%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.06"
$ws.Range("E5").Value = "  +4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.22"
$ws.Range("E6").Value = "  +8.47%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.147.87"
$ws.Range("E8").Value = "  +5.89%  "
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  +19.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  +9.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +5.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +11.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.95"
$ws.Range("E14").Value = "  +7.33%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.680.11"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.301.63"
$ws.Range("E17").Value = "  +8.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.152.95"
$ws.Range("E19").Value = "  +6.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.73"
$ws.Range("E20").Value = "  +9.50%  "
$ws.Range("E21").Value = "  +4.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  +8.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.45"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.79"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +13.04%  "
$ws.Range("E28").Value = "  +6.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.90"
$ws.Range("E31").Value = "  +12.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.34"
$ws.Range("E32").Value = "  +6.66%  "
$ws.Range("E33").Value = "  +6.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0889"
$ws.Range("E34").Value = "  +17.12%  "
$ws.Range("E35").Value = "  +18.90%  "
$ws.Range("E36").Value = "  +7.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +22.72%  "
$ws.Range("E38").Value = "  +5.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.93"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "452.61"
$ws.Range("E40").Value = "  +14.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.79"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.945.99"
$ws.Range("E42").Value = "  +8.66%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0374"
$ws.Range("E43").Value = "  +6.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.284"
$ws.Range("E44").Value = "  +14.35%  "
$ws.Range("E45").Value = "  +7.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  +13.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.48"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.17"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.06"
$ws.Range("E51").Value = "  +8.45%  "
